# Apply price/volume/coin-listing refresh for cryptos.xlsx (GitHub Actions data pull).
# Rows 26-51 also shift up by one ranking slot and gain one new coin (EOS) at the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($RangeRef, [string]$Text, [bool]$ForceText = $false)
    $value = $Text
    if ($ForceText) {
        # Leading apostrophe keeps Excel from re-parsing a numeric-looking
        # string (e.g. "8.560") back into a float and dropping the trailing zero.
        $value = "'" + $Text
    }
    $ws.Range($RangeRef).Value = $value
}

Set-CellText "D2" "30.287.45" $false
Set-CellText "E2" "  +1.77%  " $false
Set-CellText "D3" "2.090.74" $false
Set-CellText "E3" "  -0.55%  " $false
Set-CellText "E4" "  -0.59%  " $false
Set-CellText "D5" "341.65" $true
Set-CellText "E5" "  -1.71%  " $false
Set-CellText "D6" "1.002" $true
Set-CellText "E6" "  -0.55%  " $false
Set-CellText "D7" "0.5312" $true
Set-CellText "E7" "  +2.09%  " $false
Set-CellText "D8" "0.4388" $true
Set-CellText "E8" "  -0.63%  " $false
Set-CellText "D9" "54.55" $true
Set-CellText "E9" "  +0.91%  " $false
Set-CellText "D10" "0.09352" $true
Set-CellText "E10" "  -0.41%  " $false
Set-CellText "E11" "  +0.19%  " $false
Set-CellText "D12" "24.70" $true
Set-CellText "E12" "  -0.90%  " $false
Set-CellText "D13" "8.560" $true
Set-CellText "E13" "  +3.62%  " $false
Set-CellText "D14" "6.886" $true
Set-CellText "E14" "  +0.90%  " $false
Set-CellText "D15" "2.034.58" $false
Set-CellText "E15" "  -3.16%  " $false
Set-CellText "D16" "101.61" $true
Set-CellText "E16" "  -1.11%  " $false
Set-CellText "D17" "0.00001161" $true
Set-CellText "E17" "  +0.05%  " $false
Set-CellText "E18" "  -0.57%  " $false
Set-CellText "D19" "21.12" $true
Set-CellText "E19" "  -0.06%  " $false
Set-CellText "E20" "  +0.66%  " $false
Set-CellText "D21" "6.351" $true
Set-CellText "E21" "  +1.00%  " $false
Set-CellText "E22" "  -0.59%  " $false
Set-CellText "D23" "30.272.52" $false
Set-CellText "E23" "  +1.59%  " $false
Set-CellText "D24" "12.49" $true
Set-CellText "E24" "  -1.10%  " $false
Set-CellText "D25" "2.309" $true
Set-CellText "E25" "  -0.33%  " $false
Set-CellText "B26" "EthereumClassic" $false
Set-CellText "C26" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc" $false
Set-CellText "D26" "21.82" $true
Set-CellText "E26" "  -0.86%  " $false
Set-CellText "B27" "InternetComputer(DFINITY)" $false
Set-CellText "C27" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp" $false
Set-CellText "D27" "6.871" $true
Set-CellText "E27" "  +7.99%  " $false
Set-CellText "B28" "Monero" $false
Set-CellText "C28" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" $false
Set-CellText "D28" "162.79" $true
Set-CellText "E28" "  +0.22%  " $false
Set-CellText "B29" "LidoDAOToken" $false
Set-CellText "C29" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo" $false
Set-CellText "D29" "2.498" $true
Set-CellText "E29" "  -0.97%  " $false
Set-CellText "B30" "BitcoinCash" $false
Set-CellText "C30" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch" $false
Set-CellText "D30" "133.62" $true
Set-CellText "E30" "  -0.04%  " $false
Set-CellText "B31" "ImmutableX" $false
Set-CellText "C31" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" $false
Set-CellText "D31" "1.135" $true
Set-CellText "E31" "  -0.30%  " $false
Set-CellText "B32" "ARBITRUM" $false
Set-CellText "C32" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb" $false
Set-CellText "D32" "1.668" $true
Set-CellText "E32" "  -3.65%  " $false
Set-CellText "B33" "Stellar" $false
Set-CellText "C33" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" $false
Set-CellText "D33" "0.1051" $true
Set-CellText "E33" "  -0.29%  " $false
Set-CellText "B34" "Filecoin" $false
Set-CellText "C34" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil" $false
Set-CellText "D34" "6.273" $true
Set-CellText "E34" "  +0.84%  " $false
Set-CellText "B35" "HuobiToken" $false
Set-CellText "C35" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht" $false
Set-CellText "D35" "3.911" $true
Set-CellText "E35" "  -0.97%  " $false
Set-CellText "B36" "FraxShare" $false
Set-CellText "C36" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs" $false
Set-CellText "D36" "10.12" $true
Set-CellText "E36" "  -3.52%  " $false
Set-CellText "B37" "VeChain" $false
Set-CellText "C37" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" $false
Set-CellText "D37" "0.02623" $true
Set-CellText "E37" "  +1.40%  " $false
Set-CellText "B38" "Hedera" $false
Set-CellText "C38" "https://coinranking.com/coin/jad286TjB+hedera-hbar" $false
Set-CellText "D38" "0.06759" $true
Set-CellText "E38" "  +0.09%  " $false
Set-CellText "B39" "Aptos" $false
Set-CellText "C39" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt" $false
Set-CellText "D39" "12.59" $true
Set-CellText "E39" "  -0.03%  " $false
Set-CellText "B40" "TheSandbox" $false
Set-CellText "C40" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand" $false
Set-CellText "D40" "0.6956" $true
Set-CellText "E40" "  -0.62%  " $false
Set-CellText "B41" "TrustWalletToken" $false
Set-CellText "C41" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" $false
Set-CellText "D41" "1.341" $true
Set-CellText "E41" "  +0.49%  " $false
Set-CellText "B42" "Algorand" $false
Set-CellText "C42" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo" $false
Set-CellText "D42" "0.2210" $true
Set-CellText "E42" "  -0.74%  " $false
Set-CellText "B43" "Decentraland" $false
Set-CellText "C43" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana" $false
Set-CellText "D43" "0.6754" $true
Set-CellText "E43" "  -1.09%  " $false
Set-CellText "B44" "NEARProtocol" $false
Set-CellText "C44" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near" $false
Set-CellText "D44" "2.384" $true
Set-CellText "E44" "  +1.11%  " $false
Set-CellText "B45" "EnergySwap" $false
Set-CellText "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" $false
Set-CellText "D45" "14.25" $true
Set-CellText "E45" "  -1.28%  " $false
Set-CellText "B46" "Frax" $false
Set-CellText "C46" "https://coinranking.com/coin/KfWtaeV1W+frax-frax" $false
Set-CellText "D46" "1.001" $true
Set-CellText "E46" "  -0.55%  " $false
Set-CellText "B47" "WEMIXTOKEN" $false
Set-CellText "C47" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix" $false
Set-CellText "D47" "1.285" $true
Set-CellText "E47" "  +5.33%  " $false
Set-CellText "B48" "PancakeSwap" $false
Set-CellText "C48" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake" $false
Set-CellText "D48" "3.630" $true
Set-CellText "E48" "  -0.14%  " $false
Set-CellText "B49" "BabyDogeCoin" $false
Set-CellText "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge" $false
Set-CellText "D49" "0.00000000349" $true
Set-CellText "E49" "  -1.97%  " $false
Set-CellText "B50" "ThetaToken" $false
Set-CellText "C50" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta" $false
Set-CellText "D50" "1.205" $true
Set-CellText "E50" "  +3.10%  " $false
Set-CellText "B51" "EOS" $false
Set-CellText "C51" "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos" $false
Set-CellText "D51" "1.212" $true
Set-CellText "E51" "  -0.85%  " $false
